$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain text (not auto-converted to a
# number/date by Excel) while keeping the cell style unchanged, exactly like
# the original inlineStr cells in the source workbook.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2
$ws.Cells.Item(2, 4).Value = "41.232.59"
$ws.Cells.Item(2, 5).Value = "  -1.78%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.178.23"
$ws.Cells.Item(3, 5).Value = "  -2.29%  "

# Row 4
Set-TextValue 4 4 "1.00"
$ws.Cells.Item(4, 5).Value = "  -0.08%  "

# Row 5
Set-TextValue 5 4 "251.06"
$ws.Cells.Item(5, 5).Value = "  +0.29%  "

# Row 6
Set-TextValue 6 4 "0.613"
$ws.Cells.Item(6, 5).Value = "  -2.72%  "

# Row 7
Set-TextValue 7 4 "66.63"
$ws.Cells.Item(7, 5).Value = "  -7.94%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.00%  "

# Row 9
Set-TextValue 9 4 "0.579"
$ws.Cells.Item(9, 5).Value = "  -2.54%  "

# Row 10
Set-TextValue 10 4 "59.20"
$ws.Cells.Item(10, 5).Value = "  +1.78%  "

# Row 11
Set-TextValue 11 4 "36.50"
$ws.Cells.Item(11, 5).Value = "  -11.06%  "

# Row 12
Set-TextValue 12 4 "0.0938"
$ws.Cells.Item(12, 5).Value = "  -3.22%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -1.23%  "

# Row 14
Set-TextValue 14 4 "6.88"
$ws.Cells.Item(14, 5).Value = "  -4.00%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "2.503.96"
$ws.Cells.Item(15, 5).Value = "  -2.25%  "

# Row 16
Set-TextValue 16 4 "14.36"
$ws.Cells.Item(16, 5).Value = "  -4.27%  "

# Row 17
Set-TextValue 17 4 "0.849"
$ws.Cells.Item(17, 5).Value = "  -1.96%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "2.171.08"
$ws.Cells.Item(18, 5).Value = "  -2.52%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "41.132.99"
$ws.Cells.Item(19, 5).Value = "  -1.71%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "0.0₃0950"
$ws.Cells.Item(20, 5).Value = "  -1.74%  "

# Row 21
Set-TextValue 21 4 "71.73"
$ws.Cells.Item(21, 5).Value = "  -1.67%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  -2.16%  "

# Row 23
Set-TextValue 23 4 "230.67"
$ws.Cells.Item(23, 5).Value = "  -1.99%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -3.26%  "

# Row 25
Set-TextValue 25 4 "3.82"
$ws.Cells.Item(25, 5).Value = "  -9.37%  "

# Row 26
Set-TextValue 26 4 "11.47"
$ws.Cells.Item(26, 5).Value = "  +7.23%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +0.08%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "Toncoin"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue 29 4 "2.12"
$ws.Cells.Item(29, 5).Value = "  +0.57%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "Monero"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue 30 4 "168.38"
$ws.Cells.Item(30, 5).Value = "  -1.69%  "

# Row 31
Set-TextValue 31 4 "20.28"
$ws.Cells.Item(31, 5).Value = "  -2.35%  "

# Row 32
Set-TextValue 32 4 "0.122"
$ws.Cells.Item(32, 5).Value = "  -2.26%  "

# Row 33
Set-TextValue 33 4 "5.88"
$ws.Cells.Item(33, 5).Value = "  +5.05%  "

# Row 34
Set-TextValue 34 4 "0.0758"
$ws.Cells.Item(34, 5).Value = "  +3.41%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  -3.17%  "

# Row 36
Set-TextValue 36 4 "4.54"
$ws.Cells.Item(36, 5).Value = "  -4.38%  "

# Row 37
Set-TextValue 37 4 "3.97"
$ws.Cells.Item(37, 5).Value = "  +0.23%  "

# Row 38
Set-TextValue 38 4 "24.70"
$ws.Cells.Item(38, 5).Value = "  -5.90%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -0.01%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -2.94%  "

# Row 41
Set-TextValue 41 4 "5.34"
$ws.Cells.Item(41, 5).Value = "  +8.11%  "

# Row 42
Set-TextValue 42 4 "5.50"
$ws.Cells.Item(42, 5).Value = "  -8.10%  "

# Row 43
Set-TextValue 43 4 "11.49"
$ws.Cells.Item(43, 5).Value = "  -4.68%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -9.23%  "

# Row 45
Set-TextValue 45 4 "8.53"
$ws.Cells.Item(45, 5).Value = "  -2.77%  "

# Row 46
Set-TextValue 46 4 "0.0999"
$ws.Cells.Item(46, 5).Value = "  -1.78%  "

# Row 47
Set-TextValue 47 4 "0.190"
$ws.Cells.Item(47, 5).Value = "  -6.34%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -0.22%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -1.65%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  -4.15%  "

# Row 51
Set-TextValue 51 4 "2.73"
$ws.Cells.Item(51, 5).Value = "  +1.03%  "
